# Update the cryptocurrency price/volume snapshot in place.
# Values are written with a leading apostrophe via Value2 so that numeric-
# looking strings (prices with dot-grouped thousands, tiny decimals, etc.)
# and percentage strings stay literal text instead of being reinterpreted
# by Excel as numbers/dates.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'27.325.66"
$ws.Range("E2").Value2 = "'  -3.27%  "

$ws.Range("D3").Value2 = "'1.853.58"
$ws.Range("E3").Value2 = "'  -3.40%  "

$ws.Range("E4").Value2 = "'  +0.21%  "

$ws.Range("D5").Value2 = "'328.06"
$ws.Range("E5").Value2 = "'  -0.13%  "

$ws.Range("E6").Value2 = "'  +0.23%  "

$ws.Range("D7").Value2 = "'0.4594"
$ws.Range("E7").Value2 = "'  -2.10%  "

$ws.Range("D8").Value2 = "'0.3932"
$ws.Range("E8").Value2 = "'  -2.55%  "

$ws.Range("D9").Value2 = "'46.85"
$ws.Range("E9").Value2 = "'  -11.65%  "

$ws.Range("D10").Value2 = "'0.07924"
$ws.Range("E10").Value2 = "'  -5.84%  "

$ws.Range("E11").Value2 = "'  -3.50%  "

$ws.Range("D12").Value2 = "'21.44"
$ws.Range("E12").Value2 = "'  -3.35%  "

$ws.Range("D13").Value2 = "'1.852.51"
$ws.Range("E13").Value2 = "'  -3.02%  "

$ws.Range("D14").Value2 = "'5.915"
$ws.Range("E14").Value2 = "'  -2.74%  "

$ws.Range("D15").Value2 = "'7.134"
$ws.Range("E15").Value2 = "'  -4.36%  "

$ws.Range("D16").Value2 = "'1.004"
$ws.Range("E16").Value2 = "'  +0.18%  "

$ws.Range("D17").Value2 = "'86.03"
$ws.Range("E17").Value2 = "'  -4.24%  "

$ws.Range("B18").Value2 = "'TRON"
$ws.Range("C18").Value2 = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value2 = "'0.06596"
$ws.Range("E18").Value2 = "'  -0.02%  "

$ws.Range("B19").Value2 = "'ShibaInu"
$ws.Range("C19").Value2 = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value2 = "'0.00001027"
$ws.Range("E19").Value2 = "'  -3.45%  "

$ws.Range("D20").Value2 = "'17.21"
$ws.Range("E20").Value2 = "'  -4.47%  "

$ws.Range("D21").Value2 = "'1.003"
$ws.Range("E21").Value2 = "'  +0.29%  "

$ws.Range("D22").Value2 = "'5.467"
$ws.Range("E22").Value2 = "'  -4.79%  "

$ws.Range("D23").Value2 = "'27.341.47"
$ws.Range("E23").Value2 = "'  -3.24%  "

$ws.Range("D24").Value2 = "'10.87"
$ws.Range("E24").Value2 = "'  -3.36%  "

$ws.Range("D25").Value2 = "'2.293"
$ws.Range("E25").Value2 = "'  +0.75%  "

$ws.Range("D26").Value2 = "'2.073.61"
$ws.Range("E26").Value2 = "'  -3.04%  "

$ws.Range("D27").Value2 = "'153.47"
$ws.Range("E27").Value2 = "'  +0.12%  "

$ws.Range("D28").Value2 = "'20.16"
$ws.Range("E28").Value2 = "'  +0.73%  "

$ws.Range("D29").Value2 = "'2.062"
$ws.Range("E29").Value2 = "'  -2.84%  "

$ws.Range("D30").Value2 = "'5.448"
$ws.Range("E30").Value2 = "'  -5.48%  "

$ws.Range("D31").Value2 = "'121.53"
$ws.Range("E31").Value2 = "'  -1.70%  "

$ws.Range("D32").Value2 = "'0.09403"
$ws.Range("E32").Value2 = "'  -2.51%  "

$ws.Range("D33").Value2 = "'0.9473"
$ws.Range("E33").Value2 = "'  -2.95%  "

$ws.Range("E34").Value2 = "'  -1.22%  "

$ws.Range("D35").Value2 = "'3.586"
$ws.Range("E35").Value2 = "'  -1.47%  "

$ws.Range("D36").Value2 = "'5.261"
$ws.Range("E36").Value2 = "'  -5.18%  "

$ws.Range("D37").Value2 = "'0.06026"
$ws.Range("E37").Value2 = "'  -2.24%  "

$ws.Range("D38").Value2 = "'0.02229"
$ws.Range("E38").Value2 = "'  -3.09%  "

$ws.Range("D39").Value2 = "'1.207"
$ws.Range("E39").Value2 = "'  -4.69%  "

$ws.Range("E40").Value2 = "'  +0.26%  "

$ws.Range("D41").Value2 = "'8.023"
$ws.Range("E41").Value2 = "'  -9.06%  "

$ws.Range("D42").Value2 = "'0.5927"
$ws.Range("E42").Value2 = "'  -3.70%  "

$ws.Range("D43").Value2 = "'0.1887"
$ws.Range("E43").Value2 = "'  -0.94%  "

$ws.Range("D44").Value2 = "'10.17"
$ws.Range("E44").Value2 = "'  -7.96%  "

$ws.Range("D45").Value2 = "'1.283"
$ws.Range("E45").Value2 = "'  -1.51%  "

$ws.Range("D46").Value2 = "'0.5602"
$ws.Range("E46").Value2 = "'  -4.37%  "

$ws.Range("E47").Value2 = "'  -6.14%  "

$ws.Range("D48").Value2 = "'3.395"
$ws.Range("E48").Value2 = "'  -1.10%  "

$ws.Range("E49").Value2 = "'  -5.54%  "

$ws.Range("D50").Value2 = "'0.06751"
$ws.Range("E50").Value2 = "'  -2.16%  "

$ws.Range("B51").Value2 = "'Quant"
$ws.Range("C51").Value2 = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value2 = "'108.86"
$ws.Range("E51").Value2 = "'  -1.22%  "
